# Daily update at 8 AM UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the "date only" number format currently used by the last row
# (row 39) before we touch anything, so we can move it to the new last row.
$lastRowFormat = $ws.Range("A39").NumberFormat

# The previous last row (row 39) used the "date only" format.
# Now that a new row is appended after it, row 39 reverts to the
# standard "date + time" format used by all the other data rows.
$ws.Range("A39").NumberFormat = $ws.Range("A38").NumberFormat

# Append the new day's data as row 40.
$ws.Range("A40").Value = 45780
$ws.Range("B40").Value = 163
$ws.Range("C40").Value = 171
$ws.Range("D40").Value = 164

# The new last row gets the "date only" format that row 39 used to have.
$ws.Range("A40").NumberFormat = $lastRowFormat
